$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '309.12'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.32%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '11'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.51%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '11'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.215'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.38%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '11'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.86%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '11'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.640'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.07%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '11'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9148'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.29%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '11'

$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '11'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '11.22%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '11'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1824'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.57%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '11'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09176'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.41%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '11'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04224'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.03%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '11'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1052'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.09%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '11'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001251'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.33%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '11'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005812'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.70%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '11'

$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.347'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.02%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '11'

$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.317'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.50%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '11'

$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.3336'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.31%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '11'

$ws.Range("B19").Value = 'MCDex'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.420'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '11.28%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '11'

$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1382'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.22%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '11'

$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2819'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.60%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '11'

$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04026'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-1.14%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '11'

$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001265'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.72%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '11'

$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004106'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.19%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '11'

$ws.Range("B25").Value = 'NitroEx'
$ws.Range("C25").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001301'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.00%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '11'

$ws.Range("B26").Value = 'Spectre.aiUtilityToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '--'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '--%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '11'

$ws.Range("B27").Value = 'LegolasExchange'
$ws.Range("C27").Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '11'

$ws.Range("B28").Value = 'BitZToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '11'

$ws.Range("B29").Value = 'Birake'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '11'

$ws.Range("B30").Value = 'NashExchange'
$ws.Range("C30").Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '11'

$ws.Range("B31").Value = 'AAXToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '11'

$ws.Range("B32").Value = 'CenX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '11'

$ws.Range("B33").Value = 'BNIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '11'

$ws.Range("B34").Value = 'UpBots'
$ws.Range("C34").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '11'

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '11'

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '11'

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '11'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02574'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '7.71%'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '11'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05348'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '3.14%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '11'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007850'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.95%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '11'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.15%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '11'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006685'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-5.28%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '11'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.61%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '11'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008059'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.15%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '11'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3069'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.44%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '11'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006696'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.51%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '11'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.04%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '11'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1781'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '396.03%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '11'

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '11'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.04%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '11'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.04%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '11'

